# Updates crypto price/volume data in Sheet1 (columns D = Price, E = Volume(1h))
# to match the latest scrape. Price-column cells must stay plain text (some
# values look numeric, e.g. "24.40" or "1.80", and would otherwise lose their
# trailing zero / be reinterpreted as a Number if assigned directly), so the
# NumberFormat is forced to Text ("@") before each Price cell is written.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.849.99"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.376.95"
$ws.Range("E3").Value = "  -1.49%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "556.48"
$ws.Range("E5").Value = "  +1.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.61"
$ws.Range("E6").Value = "  -2.36%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.584"
$ws.Range("E8").Value = "  -1.38%  "
$ws.Range("E9").Value = "  +0.61%  "
$ws.Range("E10").Value = "  -0.36%  "
$ws.Range("E11").Value = "  +1.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.342"
$ws.Range("E12").Value = "  -3.27%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.40"
$ws.Range("E13").Value = "  -3.57%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.796.09"
$ws.Range("E14").Value = "  -1.57%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "59.802.11"
$ws.Range("E15").Value = "  +0.02%  "
$ws.Range("E16").Value = "  +0.47%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.374.14"
$ws.Range("E17").Value = "  -0.72%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.11"
$ws.Range("E18").Value = "  -2.00%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.48"
$ws.Range("E19").Value = "  +1.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "321.01"
$ws.Range("E20").Value = "  -2.58%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.64"
$ws.Range("E21").Value = "  -0.26%  "
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.11"
$ws.Range("E23").Value = "  -3.33%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.173"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("E26").Value = "  -2.51%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.38"
$ws.Range("E27").Value = "  +0.56%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.80"
$ws.Range("E28").Value = "  +1.80%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0759"
$ws.Range("E29").Value = "  -1.56%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "169.73"
$ws.Range("E30").Value = "  +0.46%  "
$ws.Range("E31").Value = "  +0.92%  "
$ws.Range("E32").Value = "  +10.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.400"
$ws.Range("E33").Value = "  -1.51%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.13"
$ws.Range("E34").Value = "  -2.63%  "
$ws.Range("E36").Value = "  +1.24%  "
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.13"
$ws.Range("E38").Value = "  -1.30%  "
$ws.Range("E39").Value = "  -0.74%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "318.21"
$ws.Range("E40").Value = "  -0.42%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "38.56"
$ws.Range("E41").Value = "  -2.18%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "144.69"
$ws.Range("E42").Value = "  +3.48%  "
$ws.Range("E43").Value = "  -3.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0970"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.67"
$ws.Range("E45").Value = "  +0.99%  "
$ws.Range("E46").Value = "  -0.63%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.568"
$ws.Range("E47").Value = "  -1.86%  "
$ws.Range("E48").Value = "  -2.32%  "
$ws.Range("E49").Value = "  +0.25%  "
$ws.Range("E50").Value = "  -0.13%  "
$ws.Range("E51").Value = "  -1.56%  "
